# JH -- Length of code -- Calculate length of a particular code.
#
# Column D holds the "length" (duration) of each coded segment, computed
# as End (column C) minus Beginning (column B). Times in B/C are stored
# as plain text in the form H:MM:SS,d (comma as decimal separator).
# We parse them to seconds, subtract, and format the result back into
# the same H:MM:SS[,ffffff] text style used by the existing D2 sample.

function ParseTimeToSeconds($t) {
    $parts = $t.Split(":")
    $h = [int]$parts[0]
    $m = [int]$parts[1]
    $secPart = $parts[2]
    if ($secPart.Contains(",")) {
        $secParts = $secPart.Split(",")
        $s = [int]$secParts[0]
        $frac = [double]("0.{0}" -f $secParts[1])
    } else {
        $s = [int]$secPart
        $frac = 0.0
    }
    return ($h * 3600) + ($m * 60) + $s + $frac
}

function FormatSecondsAsDuration($diff) {
    $h = [int][math]::Floor($diff / 3600)
    $rem = $diff - ($h * 3600)
    $m = [int][math]::Floor($rem / 60)
    $rem2 = $rem - ($m * 60)
    $sec = [int][math]::Floor($rem2)
    $frac = $rem2 - $sec

    if ([math]::Abs($frac) -lt 0.000000001) {
        return "{0}:{1:D2}:{2:D2}" -f $h, $m, $sec
    } else {
        $fracStr = "{0:F6}" -f $frac
        $fracDigits = $fracStr.Split(".")[1]
        return "{0}:{1:D2}:{2:D2},{3}" -f $h, $m, $sec, $fracDigits
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 115

for ($r = 3; $r -le $lastRow; $r++) {
    $beginText = $ws.Cells.Item($r, 2).Value()
    $endText = $ws.Cells.Item($r, 3).Value()

    $beginSeconds = ParseTimeToSeconds($beginText)
    $endSeconds = ParseTimeToSeconds($endText)

    $lengthText = FormatSecondsAsDuration($endSeconds - $beginSeconds)

    $ws.Cells.Item($r, 4).Value = $lengthText
}
